# Auto-generated Excel COM-interop edit script
# Applies numeric cell value updates (and two cell-content clears)
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR per the target diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4546582.5
$ws.Range("I33").Value = 1622
$ws.Range("J33").Value = 12987224
$ws.Range("K33").Value = 1622
$ws.Range("L33").Value = 12987224
$ws.Range("M33").Value = -1393
$ws.Range("N33").Value = -12987682
$ws.Range("H40").Value = 1356.15
$ws.Range("I40").Value = 1351.6777
$ws.Range("J40").Value = 1396.4
$ws.Range("K40").Value = 1351.6777
$ws.Range("L40").Value = 1396.4
$ws.Range("M40").Value = -1176.6777
$ws.Range("N40").Value = -1746.4
$ws.Range("H45").Value = 6897.5
$ws.Range("I45").Value = 6800
$ws.Range("J45").Value = 6995
$ws.Range("K45").Value = 20400
$ws.Range("L45").Value = 20985
$ws.Range("M45").Value = -20208
$ws.Range("N45").Value = -21369
$ws.Range("H62").Value = 1676.25
$ws.Range("I62").Value = 1676.25
$ws.Range("K62").Value = 1676.25
$ws.Range("M62").Value = -1052.25
$ws.Range("H65").Value = 1676.25
$ws.Range("I65").Value = 1676.25
$ws.Range("K65").Value = 8381.25
$ws.Range("M65").Value = -5261.25
$ws.Range("H116").Value = 6767.3335
$ws.Range("I116").Value = 11741.4
$ws.Range("J116").Value = 2245.4546
$ws.Range("K116").Value = 11741.4
$ws.Range("L116").Value = 2245.4546
$ws.Range("M116").Value = -8299.4
$ws.Range("N116").Value = -9129.454600000001
$ws.Range("H118").Value = 469.77777
$ws.Range("I118").Value = 261.14285
$ws.Range("J118").Value = 1200
$ws.Range("K118").Value = 783.4285500000001
$ws.Range("L118").Value = 3600
$ws.Range("M118").Value = 873.5714499999999
$ws.Range("N118").Value = -6914
$ws.Range("H129").Value = 892.55316
$ws.Range("I129").Value = 497
$ws.Range("J129").Value = 901.15216
$ws.Range("K129").Value = 1491
$ws.Range("L129").Value = 2703.45648
$ws.Range("M129").Value = 3509
$ws.Range("N129").Value = -12703.45648

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 5503.5
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 10007
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 10007
$ws.Range("M16").Value = -713
$ws.Range("N16").Value = -10581
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H32").Value = 6380.7
$ws.Range("I32").Value = 5120.6
$ws.Range("J32").Value = 8480.866
$ws.Range("K32").Value = 5120.6
$ws.Range("L32").Value = 8480.866
$ws.Range("M32").Value = -4833.6
$ws.Range("N32").Value = -9054.866
$ws.Range("H45").Value = 13854.556
$ws.Range("I45").Value = 26378
$ws.Range("J45").Value = 3835.8
$ws.Range("K45").Value = 26378
$ws.Range("L45").Value = 3835.8
$ws.Range("M45").Value = -26001
$ws.Range("N45").Value = -4589.8
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 254.08333
$ws.Range("I80").Value = 62
$ws.Range("J80").Value = 318.1111
$ws.Range("K80").Value = 62
$ws.Range("L80").Value = 318.1111
$ws.Range("M80").Value = 936
$ws.Range("N80").Value = -2314.1111
$ws.Range("H83").Value = 254.08333
$ws.Range("I83").Value = 62
$ws.Range("J83").Value = 318.1111
$ws.Range("K83").Value = 310
$ws.Range("L83").Value = 1590.5555
$ws.Range("M83").Value = 4682
$ws.Range("N83").Value = -11574.5555

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6789.4346
$ws.Range("I31").Value = 1523.3889
$ws.Range("K31").Value = 1523.3889
$ws.Range("M31").Value = -1228.3889
$ws.Range("H34").Value = 6789.4346
$ws.Range("I34").Value = 1523.3889
$ws.Range("K34").Value = 1523.3889
$ws.Range("M34").Value = -1321.3889
$ws.Range("H58").Value = 990.11536
$ws.Range("I58").Value = 865.2857
$ws.Range("J58").Value = 1514.4
$ws.Range("K58").Value = 865.2857
$ws.Range("L58").Value = 1514.4
$ws.Range("M58").Value = -662.2857
$ws.Range("N58").Value = -1920.4
$ws.Range("H136").Value = 990.11536
$ws.Range("I136").Value = 865.2857
$ws.Range("J136").Value = 1514.4
$ws.Range("K136").Value = 2595.8571
$ws.Range("L136").Value = 4543.200000000001
$ws.Range("M136").Value = -45.85710000000017
$ws.Range("N136").Value = -9643.200000000001
$ws.Range("H138").Value = 24620
$ws.Range("J138").Value = 24620
$ws.Range("L138").Value = 24620
$ws.Range("N138").Value = -34900
$ws.Range("H141").Value = 13163.117
$ws.Range("J141").Value = 13163.117
$ws.Range("L141").Value = 13163.117
$ws.Range("N141").Value = -23523.117

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1688
$ws.Range("J129").Value = 1851.5
$ws.Range("L129").Value = 5554.5
$ws.Range("N129").Value = -15554.5
$ws.Range("H131").Value = 42308508
$ws.Range("I131").Value = 9091385
$ws.Range("J131").Value = 66667736
$ws.Range("K131").Value = 27274155
$ws.Range("L131").Value = 200003208
$ws.Range("M131").Value = -27269115
$ws.Range("N131").Value = -200013288

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 9251
$ws.Range("J93").Value = 9251
$ws.Range("L93").Value = 9251
$ws.Range("N93").Value = -12995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6790073.5
$ws.Range("I122").Value = 7147738.5
$ws.Range("J122").Value = 5001750
$ws.Range("K122").Value = 21443215.5
$ws.Range("L122").Value = 15005250
$ws.Range("M122").Value = -21440765.5
$ws.Range("N122").Value = -15010150

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H112").Value = 58096.75
$ws.Range("J112").Value = 58096.75
$ws.Range("L112").Value = 58096.75
$ws.Range("N112").Value = -61050.75
$ws.Range("H114").Value = 42699
$ws.Range("J114").Value = 42699
$ws.Range("L114").Value = 42699
$ws.Range("N114").Value = -51377
$ws.Range("H122").Value = 2333
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H123").Value = 29883.8
$ws.Range("J123").Value = 29883.8
$ws.Range("L123").Value = 29883.8
$ws.Range("N123").Value = -39683.8
$ws.Range("H136").Value = 1127.4445
$ws.Range("I136").Value = 738.8889
$ws.Range("J136").Value = 1904.5555
$ws.Range("K136").Value = 2216.6667
$ws.Range("L136").Value = 5713.666499999999
$ws.Range("M136").Value = 333.3332999999998
$ws.Range("N136").Value = -10813.6665
